$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "SI" marker in column A for the 9-month experiment rows (19-21)
$ws.Range("A19").Value = "SI"
$ws.Range("A20").Value = "SI"
$ws.Range("A21").Value = "SI"

# Row 20 (BO step): new experiment code
$ws.Range("D20").Value = "942_HT"

# Row 21 (Modelo final, first cut): update to the "cortes hasta 15000" variant
$ws.Range("D21").Value = "992_ZZ_lightgbm_15000"
$ws.Range("C21").Value = "Entrenar el modelo final. Cortes hasta 15000"
$ws.Range("F21").Value = "exp/ZZ9412"

# Row 20 duration, added last
$ws.Range("G20").Value = "3 días"

# Remove the now-duplicate row 22 (its content was merged into row 21 above)
$ws.Rows("22:22").Delete()

# Column C is a bit wider now (auto best-fit) to accommodate the new text
$ws.Columns("C:C").ColumnWidth = 39.8

$ws.Range("A22").Select()
